$d = $word.ActiveDocument

# Update the date line
$d.Content.Find.Execute("2025-08-19 Tuesday", $true, $true, $false, $false, $false, $true, 1, $false, "2025-08-20 Wednesday", 2) | Out-Null

# Update table cells (row-major order, 20 rows x 5 columns)
$t = $d.Tables.Item(1)
if ($t.Rows.Count -ne 20 -or $t.Columns.Count -ne 5) {
    throw "Unexpected table shape: $($t.Rows.Count) rows x $($t.Columns.Count) cols"
}
$values = @(
    "91-86=",
    "76-27=",
    "27-16=",
    "48+7=",
    "99-46=",
    "11+10=",
    "81-65=",
    "10+38=",
    "91-9=",
    "26+14=",
    "40-6=",
    "18+17=",
    "65+34=",
    "86+12=",
    "18+54=",
    "87+2=",
    "12+18=",
    "85-76=",
    "85-68=",
    "92-31=",
    "38+48=",
    "10+43=",
    "58-51=",
    "19+48=",
    "29-12=",
    "0+40=",
    "27-16=",
    "96-36=",
    "5+13=",
    "35+33=",
    "8+54=",
    "8+1=",
    "50-47=",
    "51-1=",
    "81-29=",
    "30+48=",
    "1+5=",
    "97-9=",
    "46-5=",
    "87-36=",
    "42+23=",
    "38+6=",
    "59-0=",
    "77-35=",
    "35-24=",
    "47+20=",
    "70-33=",
    "48+47=",
    "60+20=",
    "79-36=",
    "65-13=",
    "48+39=",
    "46-2=",
    "52-15=",
    "45+49=",
    "25-10=",
    "85-80=",
    "83-82=",
    "27-19=",
    "93+0=",
    "95-67=",
    "93+0=",
    "40+23=",
    "43+30=",
    "94-90=",
    "5+51=",
    "16+59=",
    "38-21=",
    "73-39=",
    "43+16=",
    "35+55=",
    "86-60=",
    "13+71=",
    "78+18=",
    "65-14=",
    "75+23=",
    "80-79=",
    "93-51=",
    "88-28=",
    "63-32=",
    "38+51=",
    "34+15=",
    "83-57=",
    "74-13=",
    "13+81=",
    "57+26=",
    "73-5=",
    "7+89=",
    "39-11=",
    "57-4=",
    "42+35=",
    "87-75=",
    "10+3=",
    "59-34=",
    "60+29=",
    "94-36=",
    "38+47=",
    "88+11=",
    "0+84=",
    "74-59="
)

$rows = 20
$cols = 5
$idx = 0
for ($r = 1; $r -le $rows; $r++) {
    for ($c = 1; $c -le $cols; $c++) {
        $cell = $t.Cell($r, $c)
        $cell.Range.Text = $values[$idx]
        $idx++
    }
}

Write-Host "Done. Updated" $idx "cells."